# Apply weekly reorder of data rows (rows 2-14) on the active sheet.
# The data rows get shuffled (permuted) while keeping the same header
# row and same set of row contents - only their vertical order changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns used in the data table (A..R)
$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R")

# Snapshot current values for rows 2..14 (data rows below the header)
$snapshot = @{}
for ($r = 2; $r -le 14; $r++) {
    $rowVals = @{}
    foreach ($c in $cols) {
        $rowVals[$c] = $ws.Range("$c$r").Value2
    }
    $snapshot[$r] = $rowVals
}

# Mapping: destination row -> source row (where the content used to live)
$mapping = @{
    2  = 8
    3  = 7
    4  = 14
    5  = 12
    6  = 6
    7  = 9
    8  = 4
    9  = 13
    10 = 3
    11 = 2
    12 = 10
    13 = 5
    14 = 11
}

foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $srcVals = $snapshot[$srcRow]
    foreach ($c in $cols) {
        $ws.Range("$c$destRow").Value2 = $srcVals[$c]
    }
}
